$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-45, 48-51.
# NumberFormat is forced to text ("@") before writing the Price values so that
# Excel does not reinterpret numeric-looking strings (e.g. "335.53", "109.60")
# as actual numbers and strip formatting such as trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.405.21"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.51"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +1.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.53"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.017"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4552"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3950"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.06"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07853"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9853"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.42"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.921.81"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.866"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.019"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.022"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.33"
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06603"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001026"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.016"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.464.52"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.374"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.72"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.270"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.143.33"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.22"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.41"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.064"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.334"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.82"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9469"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09361"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.606"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.390"
$ws.Range("E35").Value = "  +4.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.255"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06032"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02214"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.168"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.153"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.016"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5751"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.10"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1809"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.246"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5432"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07213"
$ws.Range("E49").Value = "  +5.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.884"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.60"
$ws.Range("E51").Value = "  -1.05%  "

# Rows 46/47: coin order swapped (RenderToken now ranks above EnergySwap),
# update all columns (Coin, Link, Price, Volume(1h)) accordingly.
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.296"
$ws.Range("E46").Value = "  +28.79%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.93"
$ws.Range("E47").Value = "  -0.03%  "
